$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calib gases")

# --- New "calibration carrier strength mixing" table (rows 25-29) ---
# (Header cells are populated in the same left-to-right-except-totals order the
# workbook's author used, so new shared-string entries land at the same indices.)
$ws.Range("E25").Value = "source gas"
$ws.Range("F25").Value = "Calib ppm"
$ws.Range("I25").Value = "Air"
$ws.Range("H25").Value = "CO2"
$ws.Range("G25").Value = "Q total"
$ws.Range("J25").Value = "Trace"

$ws.Range("E26").Value = 1000
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 50

$ws.Range("E27").Value = 1000
$ws.Range("F27").Value = 250
$ws.Range("G27").Value = 50

$ws.Range("E28").Value = 1000
$ws.Range("F28").Value = 500
$ws.Range("G28").Value = 50

$ws.Range("E29").Value = 1000
$ws.Range("F29").Value = 750
$ws.Range("G29").Value = 50

$ws.Range("H26:H29").Formula = "=G26*F26/E26"
$ws.Range("I26:I29").Formula = "=G26*(E26-F26)/E26"

# --- Re-enter the N/P "Fixing the Carrier Flowrate" ratio formulas as one
#     fill so Excel records them as a shared formula group spanning the
#     full N11:N17 / P11:P17 block (previously each row had its own copy).
$ws.Range("N11:N17").Formula = "=M11/A$3"
$ws.Range("P11:P17").Formula = "=O11/B$3"

# Column F is a bit wider so the new header text fits comfortably.
$ws.Columns.Item(6).ColumnWidth = 11.71

# Move the active selection to reflect where the user ended up after entering the table.
[void]$ws.Range("J31").Select()
